$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty row 2 marker (engine drops truly empty rows automatically on save,
# but make sure no stale content exists there).
$ws.Rows.Item(2).ClearContents()

# Column A holds ISO-style date text ("2024-05-23"); force Text format so Excel keeps it
# as a literal string instead of auto-converting it to a date serial value.
$ws.Range("A223:A237").NumberFormat = "@"

# Row 223
$ws.Range("A223").Value = "2024-05-23"
$ws.Range("B223").Value = "12:44:53"
$ws.Range("C223").Value = "-"
$ws.Range("D223").Value = "Etiquetadora"
$ws.Range("E223").Value = "-"
$ws.Range("F223").Value = "-"
$ws.Range("G223").Value = "-"
$ws.Range("H223").Value = "12:47:11"
$ws.Range("I223").Value = "0:02:18"

# Row 224
$ws.Range("A224").Value = "2024-05-23"
$ws.Range("B224").Value = "12:44:57"
$ws.Range("C224").Value = "-"
$ws.Range("D224").Value = "Cámara no detecta Top cover"
$ws.Range("E224").Value = "-"
$ws.Range("F224").Value = "-"
$ws.Range("G224").Value = "-"
$ws.Range("H224").Value = "12:47:12"
$ws.Range("I224").Value = "0:02:15"

# Row 225
$ws.Range("A225").Value = "2024-05-23"
$ws.Range("B225").Value = "12:45:02"
$ws.Range("C225").Value = "-"
$ws.Range("D225").Value = "Power atascado en prensa, cuesta sacar"
$ws.Range("E225").Value = "-"
$ws.Range("F225").Value = "-"
$ws.Range("G225").Value = "-"
$ws.Range("H225").Value = "12:47:10"
$ws.Range("I225").Value = "0:02:08"

# Row 226
$ws.Range("A226").Value = "2024-05-23"
$ws.Range("B226").Value = "12:45:05"
$ws.Range("C226").Value = "-"
$ws.Range("D226").Value = "No detecta presencia power CP"
$ws.Range("E226").Value = "-"
$ws.Range("F226").Value = "-"
$ws.Range("G226").Value = "-"
$ws.Range("H226").Value = "12:47:10"
$ws.Range("I226").Value = "0:02:05"

# Row 227
$ws.Range("A227").Value = "2024-05-23"
$ws.Range("B227").Value = "12:45:15"
$ws.Range("C227").Value = "-"
$ws.Range("D227").Value = "Etiquetadora"
$ws.Range("E227").Value = "-"
$ws.Range("F227").Value = "-"
$ws.Range("G227").Value = "-"
$ws.Range("H227").Value = "12:47:09"
$ws.Range("I227").Value = "0:01:54"

# Row 228
$ws.Range("A228").Value = "2024-05-23"
$ws.Range("B228").Value = "12:45:59"
$ws.Range("C228").Value = "-"
$ws.Range("D228").Value = "Etiquetadora"
$ws.Range("E228").Value = "-"
$ws.Range("F228").Value = "-"
$ws.Range("G228").Value = "-"
$ws.Range("H228").Value = "12:47:08"
$ws.Range("I228").Value = "0:01:09"

# Row 229
$ws.Range("A229").Value = "2024-05-23"
$ws.Range("B229").Value = "12:46:01"
$ws.Range("C229").Value = "-"
$ws.Range("D229").Value = "Etiquetadora"
$ws.Range("E229").Value = "-"
$ws.Range("F229").Value = "-"
$ws.Range("G229").Value = "-"
$ws.Range("H229").Value = "12:47:07"
$ws.Range("I229").Value = "0:01:06"

# Row 230
$ws.Range("A230").Value = "2024-05-23"
$ws.Range("B230").Value = "12:47:04"
$ws.Range("C230").Value = "-"
$ws.Range("D230").Value = "Etiquetadora"
$ws.Range("E230").Value = "-"
$ws.Range("F230").Value = "-"
$ws.Range("G230").Value = "-"
$ws.Range("H230").Value = "12:47:05"
$ws.Range("I230").Value = "0:00:01"

# Row 231
$ws.Range("A231").Value = "2024-05-23"
$ws.Range("B231").Value = "12:50:44"
$ws.Range("C231").Value = "-"
$ws.Range("D231").Value = "Cámara no detecta foam derecho"
$ws.Range("E231").Value = "-"
$ws.Range("F231").Value = "-"
$ws.Range("G231").Value = "-"
$ws.Range("H231").Value = "12:51:15"
$ws.Range("I231").Value = "0:00:31"

# Row 232
$ws.Range("A232").Value = "2024-05-23"
$ws.Range("B232").Value = "12:54:54"
$ws.Range("C232").Value = "-"
$ws.Range("D232").Value = "Cámara no detecta Pcb"
$ws.Range("E232").Value = "-"
$ws.Range("F232").Value = "-"
$ws.Range("G232").Value = "-"
$ws.Range("H232").Value = "12:58:51"
$ws.Range("I232").Value = "0:03:57"

# Row 233
$ws.Range("A233").Value = "2024-05-23"
$ws.Range("B233").Value = "12:55:18"
$ws.Range("C233").Value = "-"
$ws.Range("D233").Value = "AOI (fallo etiqueta)"
$ws.Range("E233").Value = "-"
$ws.Range("F233").Value = "-"
$ws.Range("G233").Value = "-"
$ws.Range("H233").Value = "12:58:53"
$ws.Range("I233").Value = "0:03:35"

# Row 234
$ws.Range("A234").Value = "2024-05-23"
$ws.Range("B234").Value = "12:55:21"
$ws.Range("C234").Value = "-"
$ws.Range("D234").Value = "Etiquetadora"
$ws.Range("E234").Value = "-"
$ws.Range("F234").Value = "-"
$ws.Range("G234").Value = "-"
$ws.Range("H234").Value = "12:58:54"
$ws.Range("I234").Value = "0:03:33"

# Row 235
$ws.Range("A235").Value = "2024-05-23"
$ws.Range("B235").Value = "13:00:11"
$ws.Range("C235").Value = "-"
$ws.Range("D235").Value = "Robot no coge busbar"
$ws.Range("E235").Value = "-"
$ws.Range("F235").Value = "-"
$ws.Range("G235").Value = "-"
$ws.Range("H235").Value = "13:00:15"
$ws.Range("I235").Value = "0:00:04"

# Row 236
$ws.Range("A236").Value = "2024-05-23"
$ws.Range("B236").Value = "13:00:22"
$ws.Range("C236").Value = "-"
$ws.Range("D236").Value = "Robot no coge busbar"
$ws.Range("E236").Value = "-"
$ws.Range("F236").Value = "-"
$ws.Range("G236").Value = "-"
$ws.Range("H236").Value = "13:00:23"
$ws.Range("I236").Value = "0:00:01"

# Row 237
$ws.Range("A237").Value = "2024-05-23"
$ws.Range("B237").Value = "13:01:41"
$ws.Range("C237").Value = "-"
$ws.Range("D237").Value = "Etiquetadora"
$ws.Range("E237").Value = "-"
$ws.Range("F237").Value = "-"
$ws.Range("G237").Value = "-"
